$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-53 down to 11-54
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new data record
$ws.Cells.Item(10,1).Value = 4
$ws.Cells.Item(10,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(10,3).Value = "Los Lagos"
$ws.Cells.Item(10,4).Value = 44575
$ws.Cells.Item(10,5).Value = 10
$ws.Cells.Item(10,6).Value = 100112031
$ws.Cells.Item(10,7).Value = "Poroto verde"
$ws.Cells.Item(10,8).Value = "Sin especificar"
$ws.Cells.Item(10,9).Value = "Primera"
$ws.Cells.Item(10,10).Value = 60
$ws.Cells.Item(10,11).Value = 35000
$ws.Cells.Item(10,12).Value = 35000
$ws.Cells.Item(10,13).Value = 35000
$ws.Cells.Item(10,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(10,15).Value = "Provincia de Limarí"
$ws.Cells.Item(10,16).Value = 1400
$ws.Cells.Item(10,17).Value = 25
$ws.Cells.Item(10,18).Value = "Hortaliza"
